# Scheduled market-data refresh: updates the Leve profit-calculation columns
# (H:currentAveragePrice, I:currentAveragePriceNQ, J:currentAveragePriceHQ,
#  K:LevePriceNQ, L:LevePriceHQ, M:LeveProfitNQ, N:LeveProfitHQ) with refreshed
# Universalis market-board averages, per job sheet (ALC/ARM/BSM/CRP/CUL/GSM/LTW).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 62: The Mustache Suits Him
$ws.Range("H62").Value = 1922
$ws.Range("I62").Value = 1939.4286
$ws.Range("J62").Value = 1800
$ws.Range("K62").Value = 1939.4286
$ws.Range("L62").Value = 1800
$ws.Range("M62").Value = -1315.4286
$ws.Range("N62").Value = -3048

# Row 65: Forgery of Convenience (L)
$ws.Range("H65").Value = 1922
$ws.Range("I65").Value = 1939.4286
$ws.Range("J65").Value = 1800
$ws.Range("K65").Value = 9697.143
$ws.Range("L65").Value = 9000
$ws.Range("M65").Value = -6577.143
$ws.Range("N65").Value = -15240

# Row 76: Warding Off Temptation
$ws.Range("H76").Value = 18153.666
$ws.Range("J76").Value = 7727
$ws.Range("L76").Value = 7727
$ws.Range("N76").Value = -8357

# Row 79: The Garden of Arcane Delights (L)
$ws.Range("H79").Value = 18153.666
$ws.Range("J79").Value = 7727
$ws.Range("L79").Value = 7727
$ws.Range("N79").Value = -9911

# Row 80: Cleansing the Wicked Humours
$ws.Range("H80").Value = 1574.2632
$ws.Range("I80").Value = 1362.6154
$ws.Range("J80").Value = 2032.8334
$ws.Range("K80").Value = 4087.8462
$ws.Range("L80").Value = 6098.5002
$ws.Range("M80").Value = -3089.8462
$ws.Range("N80").Value = -8094.5002

# Row 83: Washing Away the Sins (L)
$ws.Range("H83").Value = 1574.2632
$ws.Range("I83").Value = 1362.6154
$ws.Range("J83").Value = 2032.8334
$ws.Range("K83").Value = 12263.5386
$ws.Range("L83").Value = 18295.5006
$ws.Range("M83").Value = -7271.5386
$ws.Range("N83").Value = -28279.5006

# Row 88: The Grave of Hemlock Groves
$ws.Range("H88").Value = 4242.857
$ws.Range("I88").Value = 2903
$ws.Range("J88").Value = 4466.1665
$ws.Range("K88").Value = 2903
$ws.Range("L88").Value = 4466.1665
$ws.Range("M88").Value = -2497
$ws.Range("N88").Value = -5278.1665

# Row 91: Dappling the Highlands (L)
$ws.Range("H91").Value = 4242.857
$ws.Range("I91").Value = 2903
$ws.Range("J91").Value = 4466.1665
$ws.Range("K91").Value = 2903
$ws.Range("L91").Value = 4466.1665
$ws.Range("M91").Value = -1499
$ws.Range("N91").Value = -7274.1665

# Row 113: Amaro Kart
$ws.Range("H113").Value = 4797.2856
$ws.Range("J113").Value = 4610.7144
$ws.Range("L113").Value = 4610.7144
$ws.Range("N113").Value = -11118.7144

# Row 118: Crafty Concoctions
$ws.Range("H118").Value = 1015.5
$ws.Range("I118").Value = 1015.5
$ws.Range("K118").Value = 3046.5
$ws.Range("M118").Value = -1389.5

# Row 132: Fast-forwarding Flora
$ws.Range("H132").Value = 2417.224
$ws.Range("I132").Value = 2555.54
$ws.Range("K132").Value = 7666.62
$ws.Range("M132").Value = -5136.62

# Row 137: Cutting Edge of Culinary Quality
$ws.Range("H137").Value = 13819.5
$ws.Range("I137").Value = 15201
$ws.Range("K137").Value = 45603
$ws.Range("M137").Value = -43053

# Row 138: All-night Crafting
$ws.Range("H138").Value = 17859720
$ws.Range("I138").Value = 1027.28
$ws.Range("K138").Value = 3081.84
$ws.Range("M138").Value = 2058.16

# Row 141: Remedy for Reason
$ws.Range("H141").Value = 1658.6316
$ws.Range("I141").Value = 1387.4375
$ws.Range("J141").Value = 3105
$ws.Range("K141").Value = 4162.3125
$ws.Range("L141").Value = 9315
$ws.Range("M141").Value = 1017.6875
$ws.Range("N141").Value = -19675

$ws = $wb.Worksheets.Item("ARM")
# Row 32: Ingot We Trust
$ws.Range("H32").Value = 3016.718
$ws.Range("I32").Value = 3009.3867
$ws.Range("K32").Value = 3009.3867
$ws.Range("M32").Value = -2722.3867

# Row 61: Dealing with the Tough Stuff
$ws.Range("H61").Value = 3096.1936
$ws.Range("I61").Value = 2877.842
$ws.Range("J61").Value = 3441.9167
$ws.Range("K61").Value = 2877.842
$ws.Range("L61").Value = 3441.9167
$ws.Range("M61").Value = -2665.842
$ws.Range("N61").Value = -3865.9167

# Row 88: The Mast Chance
$ws.Range("H88").Value = 17565.572
$ws.Range("I88").Value = 34484
$ws.Range("J88").Value = 4876.75
$ws.Range("K88").Value = 34484
$ws.Range("L88").Value = 4876.75
$ws.Range("M88").Value = -34078
$ws.Range("N88").Value = -5688.75

# Row 91: The Rose and the Riveter (L)
$ws.Range("H91").Value = 17565.572
$ws.Range("I91").Value = 34484
$ws.Range("J91").Value = 4876.75
$ws.Range("K91").Value = 34484
$ws.Range("L91").Value = 4876.75
$ws.Range("M91").Value = -33080
$ws.Range("N91").Value = -7684.75

# Row 132: Don't Bore Me, Ore Me
$ws.Range("H132").Value = 66713.95
$ws.Range("I132").Value = 6162.952
$ws.Range("J132").Value = 702499.5
$ws.Range("K132").Value = 18488.856
$ws.Range("L132").Value = 2107498.5
$ws.Range("M132").Value = -15958.856
$ws.Range("N132").Value = -2112558.5

# Row 136: Metal with Mettle
$ws.Range("H136").Value = 3096.1936
$ws.Range("I136").Value = 2877.842
$ws.Range("J136").Value = 3441.9167
$ws.Range("K136").Value = 8633.526
$ws.Range("L136").Value = 10325.7501
$ws.Range("M136").Value = -6083.526
$ws.Range("N136").Value = -15425.7501

$ws = $wb.Worksheets.Item("BSM")
# Row 107: The Gold Experience
$ws.Range("H107").Value = 3264.389
$ws.Range("I107").Value = 3051
$ws.Range("K107").Value = 3051
$ws.Range("M107").Value = -1131

# Row 134: Ruthenium Supremium
$ws.Range("H134").Value = 1050.7778
$ws.Range("I134").Value = 994.94116
$ws.Range("K134").Value = 2984.82348
$ws.Range("M134").Value = -449.82348

$ws = $wb.Worksheets.Item("CRP")
# Row 86: Birch, Please
$ws.Range("H86").Value = 14997
$ws.Range("J86").Value = 14997
$ws.Range("L86").Value = 14997
$ws.Range("N86").Value = -17243

# Row 89: Built This City on Blocks and Soul (L)
$ws.Range("H89").Value = 14997
$ws.Range("J89").Value = 14997
$ws.Range("L89").Value = 74985
$ws.Range("N89").Value = -86217

# Row 99: O Pine
$ws.Range("H99").Value = 5992.909
$ws.Range("I99").Value = 5659.778
$ws.Range("J99").Value = 7492
$ws.Range("K99").Value = 5659.778
$ws.Range("L99").Value = 7492
$ws.Range("M99").Value = -4161.778
$ws.Range("N99").Value = -10488

# Row 107: Built to Last
$ws.Range("H107").Value = 4980.25
$ws.Range("I107").Value = 3999
$ws.Range("J107").Value = 5307.3335
$ws.Range("K107").Value = 3999
$ws.Range("L107").Value = 5307.3335
$ws.Range("M107").Value = -2079
$ws.Range("N107").Value = -9147.333500000001

# Row 126: A Better Conductor
$ws.Range("H126").Value = 5992.909
$ws.Range("I126").Value = 5659.778
$ws.Range("J126").Value = 7492
$ws.Range("K126").Value = 16979.334
$ws.Range("L126").Value = 22476
$ws.Range("M126").Value = -14509.334
$ws.Range("N126").Value = -27416

$ws = $wb.Worksheets.Item("CUL")
# Row 4: In Hot Water
$ws.Range("H4").Value = 67488650
$ws.Range("I4").Value = 113080824
$ws.Range("K4").Value = 339242472
$ws.Range("M4").Value = -339242360

# Row 38: Pretty as a Picture
$ws.Range("H38").Value = 171.4
$ws.Range("I38").Value = 101.5
$ws.Range("J38").Value = 218
$ws.Range("K38").Value = 304.5
$ws.Range("L38").Value = 654
$ws.Range("M38").Value = 42.5
$ws.Range("N38").Value = -1348

# Row 88: Don't Let It Fall Apart
$ws.Range("H88").Value = 20000
$ws.Range("I88").Value = 0
$ws.Range("J88").Value = 20000
$ws.Range("K88").Value = 0
$ws.Range("L88").Value = 60000
$ws.Range("M88").ClearContents()
$ws.Range("N88").Value = -60856

# Row 91: Better Come Back with a Sandwich (L)
$ws.Range("H91").Value = 20000
$ws.Range("I91").Value = 0
$ws.Range("J91").Value = 20000
$ws.Range("K91").Value = 0
$ws.Range("L91").Value = 60000
$ws.Range("M91").ClearContents()
$ws.Range("N91").Value = -62964

$ws = $wb.Worksheets.Item("GSM")
# Row 102: Put the Metal to the Peddle
$ws.Range("H102").Value = 4023.4062
$ws.Range("I102").Value = 2098.1738
$ws.Range("J102").Value = 8943.444
$ws.Range("K102").Value = 2098.1738
$ws.Range("L102").Value = 8943.444
$ws.Range("M102").Value = -476.1738
$ws.Range("N102").Value = -12187.444

# Row 113: Copious Crystal Cannons
$ws.Range("H113").Value = 2910.389
$ws.Range("I113").Value = 2078.3333
$ws.Range("J113").Value = 4574.5
$ws.Range("K113").Value = 2078.3333
$ws.Range("L113").Value = 4574.5
$ws.Range("M113").Value = 91.66670000000022
$ws.Range("N113").Value = -8914.5

$ws = $wb.Worksheets.Item("LTW")
# Row 16: Saddle Sore
$ws.Range("H16").Value = 900.5
$ws.Range("J16").Value = 1383.3334
$ws.Range("L16").Value = 1383.3334
$ws.Range("N16").Value = -1723.3334

# Row 133: The Perfect Accessory
$ws.Range("H133").Value = 66286
$ws.Range("J133").Value = 66286
$ws.Range("L133").Value = 66286
$ws.Range("N133").Value = -71346
